$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values following a repull of data
$ws.Range("F2").Value = -9
$ws.Range("F4").Value = -4
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 5
$ws.Range("F14").Value = 5
$ws.Range("F19").Value = -2
$ws.Range("F29").Value = 9
$ws.Range("F32").Value = -4
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 1
